$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style/format of the existing header cell (G1, "sum") onto the
# new header cell H1 so it matches the bold/bordered/centered header style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New header: "Save"
$ws.Range("H1").Value = "Save"

# New data value for row 2, column H
$ws.Range("H2").Value = 0
